# Apply the author's edits:
#  1. Delete the slide "Component 2 (Trello screenshot)" (sldId=266),
#     which sits at position 14 in the deck (it is the third slide
#     sharing that same title text, so it must be targeted by its
#     unique slide Id rather than by title text).
#  2. On the "Relevant Implications" slide, append a trailing space to
#     the last run of the body text so it reads
#     "... to learn how to do this. ".

$p = $ppt.ActivePresentation

# --- 1. Delete the duplicate "Component 2 (Trello screenshot)" slide -----
# (SlideID 266 is the unique, stable identifier for this slide; three
# slides in the deck share the same title text, so match on SlideID.)
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 266) {
        $slide.Delete()
    }
}

# --- 2. Fix trailing text on the "Relevant Implications" slide -----------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 257) {
        for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
            $shp = $slide.Shapes.Item($j)
            if ($shp.HasTextFrame) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "Explain the relevant implications here.  Please watch this video to learn how to do this.") {
                    $runCount = $tr.Runs().Count
                    $lastRun = $tr.Runs($runCount)
                    if ($lastRun.Text -eq " to learn how to do this.") {
                        $lastRun.Text = " to learn how to do this. "
                    }
                }
            }
        }
    }
}
